# Added ifoCAST full series evaluation: the rolling error-table window
# shifted forward by one quarter. Each row's B:G values become the values
# that used to belong to the next row, and the final row (11) receives the
# newly computed statistics for the additional quarter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, -0.04706953193124818, 0.7263783780799662, 0.9076748761767716, 0.9527197259303344, 0.9791433241278349, 18),
    @(3, 0.008980548322284209, 0.7068694245681296, 0.7853984072668317, 0.8862270630413133, 0.9134550436843235, 17),
    @(4, 0.130843318123659,    0.6351380739751387, 0.6444494568402878, 0.8027760938395511, 0.8180167776226784, 16),
    @(5, 0.07882986512450849,  0.6136611455857488, 0.5739360866201003, 0.7575856958919567, 0.779918909922613,  15),
    @(6, 0.1364943803647256,   0.5668254508253494, 0.6024723940093601, 0.7761909520275021, 0.7929392009836149, 14),
    @(7, 0.1194123702819532,   0.6296987908481204, 0.6261569364748241, 0.7913007876116541, 0.8141800113140473, 13),
    @(8, 0.2385215453226533,   0.5677170375037683, 0.505199978683133,  0.7107742107611481, 0.6993301608877609, 12),
    @(9, 0.1002722471476047,   0.3965624711381974, 0.2295019047709758, 0.4790635706991044, 0.49131671999357,   11),
    @(10, 0.1714703583411825,  0.390033206010506,  0.2698760920897213, 0.519495998146012,  0.5169074796862885, 10),
    @(11, 0.1087828097550494,  0.3930655897532851, 0.2412922028369822, 0.4912150270879161, 0.5080756007043895, 9)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
